$d = $word.ActiveDocument

# --- Locate the "SMARTREWARDS FAQ'S" heading paragraph robustly via Find ---
$findRange = $d.Content
$findRange.Find.Execute("SMARTREWARDS FAQ") | Out-Null
if (-not $findRange.Find.Found) {
    throw "Could not find 'SMARTREWARDS FAQ' heading"
}
$headStart = $findRange.Start

# Range.Paragraphs indexing isn't reliable in this host, so resolve the
# paragraph index by scanning Document.Paragraphs for the one containing
# the Find hit.
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $headStart -and $p.Range.End -gt $headStart) {
        $headingIndex = $i
        break
    }
}
if ($headingIndex -eq -1) {
    throw "Could not resolve heading paragraph index"
}

# The trailing empty paragraph at the very end of the body must be kept,
# so the deletion stops right before it.
$totalParas = $d.Paragraphs.Count
$lastParaIndex = $totalParas

$startDelete = $d.Paragraphs.Item($headingIndex + 1).Range.Start
$endDelete = $d.Paragraphs.Item($lastParaIndex).Range.Start

$deleteRange = $d.Range($startDelete, $endDelete)
$deleteRange.Delete()
